$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.833.49'
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("D3").Value = '1.838.02'
$ws.Range("E3").Value = '  +1.25%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.619'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.94%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '39.60'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.03%  '
$ws.Range("E9").Value = '  +1.10%  '
$ws.Range("E10").Value = '  -0.10%  '
$ws.Range("E11").Value = '  -1.36%  '
$ws.Range("D12").Value = '2.105.62'
$ws.Range("E12").Value = '  +1.24%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.854.58'
$ws.Range("E13").Value = '  +2.20%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.33'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.671'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.78%  '
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").Value = '34.850.63'
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.64'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("D19").Value = '0.0₃0785'
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '240.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.16'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.14%  '
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.27'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.60%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.78'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("E27").Value = '  +2.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.39'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.37%  '
$ws.Range("E29").Value = '  -5.27%  '
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("E31").Value = '  -0.40%  '
$ws.Range("E32").Value = '  -2.73%  '
$ws.Range("E33").Value = '  -1.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.87'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.13%  '
$ws.Range("E35").Value = '  +6.96%  '
$ws.Range("E36").Value = '  +10.63%  '
$ws.Range("E37").Value = '  +1.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '91.40'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.49%  '
$ws.Range("E39").Value = '  +5.48%  '
$ws.Range("D40").Value = '1.339.95'
$ws.Range("E40").Value = '  +2.39%  '
$ws.Range("E41").Value = '  -0.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.56'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.43%  '
$ws.Range("E43").Value = '  -1.35%  '
$ws.Range("E44").Value = '  -2.71%  '
$ws.Range("E45").Value = '  -0.30%  '
$ws.Range("E46").Value = '  -0.81%  '
$ws.Range("E47").Value = '  +1.87%  '
$ws.Range("D48").Value = '2.019.37'
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0681'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.28%  '
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.27'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +14.93%  '
